$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.499.04'
$ws.Range('E2').Value = '  -5.12%  '
$ws.Range('D3').Value = '2.579.50'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '299.14'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').Value = '95.37'
$ws.Range('E6').Value = '  -3.50%  '
$ws.Range('E7').Value = '  -4.38%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.548'
$ws.Range('E9').Value = '  -4.57%  '
$ws.Range('D10').Value = '36.61'
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  -4.01%  '
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  -5.04%  '
$ws.Range('D13').Value = '2.968.19'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').Value = '2.583.35'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '0.878'
$ws.Range('E16').Value = '  -4.40%  '
$ws.Range('D17').Value = '14.17'
$ws.Range('E17').Value = '  -5.04%  '
$ws.Range('D18').Value = '43.456.27'
$ws.Range('E18').Value = '  -5.70%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '6.56'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0964'
$ws.Range('E20').Value = '  -3.93%  '
$ws.Range('D21').Value = '12.23'
$ws.Range('E21').Value = '  -5.26%  '
$ws.Range('D22').Value = '72.74'
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('D23').Value = '262.62'
$ws.Range('E23').Value = '  -3.15%  '
$ws.Range('D24').Value = '2.90'
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').Value = '29.02'
$ws.Range('E26').Value = '  -2.21%  '
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Value = '10.13'
$ws.Range('E28').Value = '  -3.75%  '
$ws.Range('E29').Value = '  -4.37%  '
$ws.Range('D30').Value = '37.34'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').Value = '5.91'
$ws.Range('E31').Value = '  -5.65%  '
$ws.Range('D32').Value = '3.57'
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').Value = '2.20'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('E34').Value = '  -1.53%  '
$ws.Range('D35').Value = '150.76'
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('D36').Value = '0.0801'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').Value = '24.34'
$ws.Range('E38').Value = '  +6.29%  '
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('D40').Value = '16.46'
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').Value = '3.42'
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').Value = '0.0311'
$ws.Range('E42').Value = '  -4.81%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '3.80'
$ws.Range('E43').Value = '  -6.49%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.070.19'
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '87.77'
$ws.Range('E46').Value = '  -5.68%  '
$ws.Range('D47').Value = '9.09'
$ws.Range('E47').Value = '  -4.37%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '1.60'
$ws.Range('E48').Value = '  +3.11%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.835.97'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = '104.97'
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').Value = '0.189'
$ws.Range('E51').Value = '  -4.77%  '
